# Insert two new columns (I and J) into the "Tablib Dataset" sheet, shifting
# the old "Export Issues"/"Location of Issue" columns to K/L, then populate
# the new columns and update the (now shifted) "Export Issues" data value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for two new columns before the old column I ("Export Issues").
$ws.Columns.Item(9).Insert()
$ws.Columns.Item(10).Insert()

# New column I: "Date of Last Update" (reuses the old "Export Issues" data,
# which is the timestamp that used to live in column I).
$ws.Cells.Item(1, 9).Value = "Date of Last Update"
$ws.Cells.Item(2, 9).Value = "2017-08-15T09:19:15.453000"

# New column J: duplicate "TicketID" header/value from column A.
$ws.Cells.Item(1, 10).Value = "TicketID"
$ws.Cells.Item(2, 10).Value = "'243"
$ws.Cells.Item(2, 10).Style = "Normal"

# Old column I ("Export Issues") is now column K - correct its data value.
$ws.Cells.Item(2, 11).Value = "Other-Not Listed"
